$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 131
$ws.Cells.Item(3, 6).Value = 335
$ws.Cells.Item(5, 6).Value = 1207
$ws.Cells.Item(6, 6).Value = 437
$ws.Cells.Item(8, 6).Value = 155
$ws.Cells.Item(11, 6).Value = 3
$ws.Cells.Item(13, 6).Value = 158
$ws.Cells.Item(14, 6).Value = 169
$ws.Cells.Item(15, 6).Value = 1427
$ws.Cells.Item(16, 6).Value = 533
$ws.Cells.Item(17, 6).Value = 212
$ws.Cells.Item(18, 6).Value = 324
$ws.Cells.Item(20, 6).Value = 764
$ws.Cells.Item(21, 6).Value = 1127
$ws.Cells.Item(23, 6).Value = 1921
$ws.Cells.Item(24, 6).Value = 2600
$ws.Cells.Item(25, 6).Value = 1362
$ws.Cells.Item(27, 6).Value = 12
$ws.Cells.Item(28, 6).Value = 317
$ws.Cells.Item(29, 6).Value = 391
$ws.Cells.Item(30, 6).Value = 1078
$ws.Cells.Item(31, 6).Value = 790
$ws.Cells.Item(32, 6).Value = 1227
$ws.Cells.Item(33, 6).Value = 141
$ws.Cells.Item(35, 6).Value = 771
$ws.Cells.Item(36, 6).Value = 521
$ws.Cells.Item(37, 6).Value = 629
$ws.Cells.Item(38, 6).Value = 804
$ws.Cells.Item(39, 6).Value = 342
$ws.Cells.Item(40, 6).Value = 223
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(15, 6).Value = 582
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(6, 6).Value = 131
$ws.Cells.Item(7, 6).Value = 335
$ws.Cells.Item(11, 6).Value = 1207
$ws.Cells.Item(12, 6).Value = 437
$ws.Cells.Item(14, 6).Value = 155
$ws.Cells.Item(19, 6).Value = 158
$ws.Cells.Item(20, 6).Value = 169
$ws.Cells.Item(21, 6).Value = 1427
$ws.Cells.Item(22, 6).Value = 533
$ws.Cells.Item(23, 6).Value = 212
$ws.Cells.Item(24, 6).Value = 324
$ws.Cells.Item(26, 6).Value = 1127
$ws.Cells.Item(27, 6).Value = 2600
$ws.Cells.Item(29, 6).Value = 1362
$ws.Cells.Item(34, 6).Value = 317
$ws.Cells.Item(35, 6).Value = 391
$ws.Cells.Item(36, 6).Value = 1079
$ws.Cells.Item(39, 6).Value = 790
$ws.Cells.Item(40, 6).Value = 1227
$ws.Cells.Item(41, 6).Value = 771
$ws.Cells.Item(42, 6).Value = 521
$ws.Cells.Item(43, 6).Value = 629
$ws.Cells.Item(44, 6).Value = 804
$ws.Cells.Item(45, 6).Value = 342
$ws.Cells.Item(48, 6).Value = 223
